$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "36d00f151300d821b1a5eb9d15927671"
$ws.Range("B11").Value = "17943777771b2c0517d6f533f3897e53"
$ws.Range("B15").Value = "19afadabbdf38755ee10461b6337a744"
$ws.Range("B17").Value = "6998a33229398d564c3b7144cc6fd9f7"
$ws.Range("B24").Value = "2840a96234e359893e99f5d6e79c69a9"
$ws.Range("B29").Value = "5b9d2fbaeb862bacd1e974541000918d"
$ws.Range("B73").Value = "6c8d95b0fb1294438b0acbe5756b5428"
$ws.Range("B121").Value = "19e084c9997c8863103f974a01f5440c"
$ws.Range("B126").Value = "6308473a7b65986a8611b07fbff53ae8"
$ws.Range("B133").Value = "1eb384b3b8327844eae80cfb91424538"
$ws.Range("B136").Value = "f914d70d0ba65ea11d4f594216794952"
$ws.Range("B159").Value = "669475437e5b5f46fa1477c92d1de3a9"
$ws.Range("B169").Value = "ec72b02bb2e6801c945a84aab7da2c60"
$ws.Range("B183").Value = "745d31d0fbddd95b7ae7e57aafa25000"
$ws.Range("B191").Value = "fe9e39708fad5f08a46bd4f0471211e1"
$ws.Range("B198").Value = "ebb887ec7a4405b3a90801d6b61c44c6"
$ws.Range("B200").Value = "88b9e284edd47815a798250350a758ae"
$ws.Range("B228").Value = "583c503acd6c5c3db76863d10d491ea3"
$ws.Range("B281").Value = "e69cf0fdfb323ff998d24bae1df1ed77"
$ws.Range("B302").Value = "df48447f3941ba614cdc5a7e425d6dbc"
$ws.Range("B339").Value = "752f3b3b0545a1405228b72a412470c0"
$ws.Range("B460").Value = "9e942eb37190eb7c6f3b434ce19c917c"
$ws.Range("B480").Value = "c6cdac5190cdc1fd4ec3b89ad6e98a7b"
$ws.Range("B500").Value = "9608a861c301660357a44eee448732c4"
$ws.Range("B501").Value = "7d7627502e3c22104bd7c10b38a3505a"
$ws.Range("B502").Value = "b112424db46bbea9a28a9febfec811ad"
$ws.Range("B515").Value = "07ec601cb5de3fa09ad746a13f659cc2"
$ws.Range("B517").Value = "2a3c894fa2b7bd03c4fdb20b3d95bded"
$ws.Range("B550").Value = "0f73398ffcfaa1146a8b7b7406f2cf0c"
$ws.Range("B572").Value = "f0f9ef3977e9c453b571322e153f27f4"
$ws.Range("B616").Value = "2575727c67ef0944d9b2fc827f747a61"
$ws.Range("B627").Value = "2aa7dbfe12e7a952c29050207f73a44a"
$ws.Range("B629").Value = "ce5a2e8e400c40d4e934c845d0bb2d5d"
$ws.Range("B649").Value = "dfedadf3e3cbfb32b6174c39a509bccf"
$ws.Range("B655").Value = "ce0506567ca36001127d8af6ea5feeec"
$ws.Range("B665").Value = "1c5fa512a7d19ff7371da9ab0e1ce20b"
$ws.Range("B680").Value = "f06d3a7290c6e23b63ee300546868111"
$ws.Range("B685").Value = "1dd606b99eb8734bb1cef36dc1c848e2"
$ws.Range("B703").Value = "c3f0f223ff89c8b2c02bfc4ed5b2ca62"
$ws.Range("B704").Value = "ef545a22b6225be18b5b0e1aafb7eed1"
$ws.Range("B715").Value = "d3d294719a1a1165b656ae2b33385fc5"
$ws.Range("B742").Value = "ff155abdfd6b2b002a160c1a57874c66"
$ws.Range("B819").Value = "d2b350bc4835484594a63b0da7301925"
$ws.Range("B830").Value = "57986762c6de8356928d6e6dbf08753d"
$ws.Range("B835").Value = "880e9b057c86aedf2f9a057bd4275d1a"
$ws.Range("B862").Value = "ac1647ddba840788552c613cc76bbbe6"
